$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Org_y" -> "Org_x"
$ws.Range("A1").Value = "Org_x"

# Append a new data row (row 19) for "Ingredion" with all zero values
$newRow = 19
$lastRow = $newRow - 1

# Match the formatting used by the existing column-A cells (bold, centered, bordered)
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "Ingredion"

for ($col = 2; $col -le 44; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 0
}
